# Update crypto price (column D) and 1h volume change % (column E) values
# reflecting refreshed market data pulled by the scraper workflow.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '42.843.92'; DForceText = $false; E = '  -1.50%  ' }
    @{ Row = 3; D = '2.298.45'; DForceText = $false; E = '  -1.87%  ' }
    @{ Row = 4; D = '1.00'; DForceText = $true; E = '  -0.07%  ' }
    @{ Row = 5; D = '300.57'; DForceText = $true; E = '  -1.44%  ' }
    @{ Row = 6; D = '96.56'; DForceText = $true; E = '  -4.62%  ' }
    @{ Row = 7; D = '0.502'; DForceText = $true; E = '  -1.19%  ' }
    @{ Row = 9; D = '0.494'; DForceText = $true; E = '  -2.50%  ' }
    @{ Row = 10; D = '33.48'; DForceText = $true; E = '  -4.80%  ' }
    @{ Row = 11; D = '0.0794'; DForceText = $true; E = '  -0.48%  ' }
    @{ Row = 12; D = '49.06'; DForceText = $true; E = '  -4.59%  ' }
    @{ Row = 13; D = $null; DForceText = $false; E = '  +2.26%  ' }
    @{ Row = 14; D = '16.83'; DForceText = $true; E = '  +6.99%  ' }
    @{ Row = 15; D = '6.77'; DForceText = $true; E = '  -0.55%  ' }
    @{ Row = 16; D = '2.650.83'; DForceText = $false; E = '  -2.23%  ' }
    @{ Row = 17; D = '2.295.52'; DForceText = $false; E = '  -2.10%  ' }
    @{ Row = 18; D = '0.798'; DForceText = $true; E = '  -1.58%  ' }
    @{ Row = 19; D = '42.708.46'; DForceText = $false; E = '  -1.62%  ' }
    @{ Row = 20; D = '0.0₃0899'; DForceText = $false; E = '  -1.11%  ' }
    @{ Row = 21; D = '11.51'; DForceText = $true; E = '  -2.78%  ' }
    @{ Row = 22; D = '5.99'; DForceText = $true; E = '  -2.17%  ' }
    @{ Row = 23; D = '67.01'; DForceText = $true; E = '  -1.70%  ' }
    @{ Row = 24; D = '236.66'; DForceText = $true; E = '  -0.48%  ' }
    @{ Row = 25; D = '1.98'; DForceText = $true; E = '  +0.02%  ' }
    @{ Row = 26; D = $null; DForceText = $false; E = '  +0.10%  ' }
    @{ Row = 27; D = '2.45'; DForceText = $true; E = '  -3.42%  ' }
    @{ Row = 28; D = '24.72'; DForceText = $true; E = '  -1.35%  ' }
    @{ Row = 29; D = $null; DForceText = $false; E = '  -1.11%  ' }
    @{ Row = 30; D = $null; DForceText = $false; E = '  +0.55%  ' }
    @{ Row = 31; D = '33.86'; DForceText = $true; E = '  -2.21%  ' }
    @{ Row = 32; D = '9.11'; DForceText = $true; E = '  -1.44%  ' }
    @{ Row = 34; D = '4.79'; DForceText = $true; E = '  +6.54%  ' }
    @{ Row = 35; D = $null; DForceText = $false; E = '  -2.36%  ' }
    @{ Row = 36; D = $null; DForceText = $false; E = '  -0.71%  ' }
    @{ Row = 37; D = '16.93'; DForceText = $true; E = '  +0.50%  ' }
    @{ Row = 38; D = '0.0694'; DForceText = $true; E = '  -1.54%  ' }
    @{ Row = 39; D = '2.81'; DForceText = $true; E = '  -3.52%  ' }
    @{ Row = 40; D = '0.100'; DForceText = $true; E = '  -2.05%  ' }
    @{ Row = 41; D = '1.74'; DForceText = $true; E = '  -4.59%  ' }
    @{ Row = 42; D = $null; DForceText = $false; E = '  -2.04%  ' }
    @{ Row = 43; D = '2.34'; DForceText = $true; E = '  -3.49%  ' }
    @{ Row = 44; D = '1.960.18'; DForceText = $false; E = '  -1.15%  ' }
    @{ Row = 45; D = '0.0281'; DForceText = $true; E = '  -1.64%  ' }
    @{ Row = 46; D = '17.65'; DForceText = $true; E = '  -4.80%  ' }
    @{ Row = 47; D = '9.78'; DForceText = $true; E = '  -2.57%  ' }
    @{ Row = 48; D = '2.83'; DForceText = $true; E = '  -3.73%  ' }
    @{ Row = 49; D = '2.520.78'; DForceText = $false; E = '  -1.87%  ' }
    @{ Row = 50; D = '52.69'; DForceText = $true; E = '  -6.58%  ' }
    @{ Row = 51; D = '4.56'; DForceText = $true; E = '  -6.11%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        if ($u.DForceText) {
            # Force text interpretation so values like '1.00' or '0.100' keep
            # their literal digits instead of being parsed as numbers, then
            # drop the temporary formatting so the cell's style is unaffected.
            $dCell.NumberFormat = "@"
            $dCell.Value = $u.D
            $dCell.ClearFormats()
        } else {
            $dCell.Value = $u.D
        }
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
